# Rename the first and third worksheets (space -> underscore), then make
# "Cell volume" (now "Cell_volume", 1st tab) the active sheet with cell
# J25 selected — matching the target view state, where the active/selected
# tab moves from "PCA input" (3rd tab, which keeps its own selection at
# E53) to "Cell volume", and that sheet's stored selection changes from
# I25 to J25.

$wb = $excel.ActiveWorkbook

$wsCellVolume = $wb.Worksheets.Item(1)   # "Cell volume"
$wsPcaInput   = $wb.Worksheets.Item(3)   # "PCA input"

$wsCellVolume.Name = "Cell_volume"
$wsPcaInput.Name   = "PCA_input"

# Make "Cell_volume" the active sheet and select J25 on it (was I25).
[void]$wsCellVolume.Activate()
[void]$wsCellVolume.Range("J25").Select()
